$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.464.01'
$ws.Range("E2").Value = '  +3.27%  '
$ws.Range("D3").Value = '1.920.46'
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Formula = '="249.16"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Formula = '="43.85"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +1.08%  '
$ws.Range("D9").Formula = '="58.56"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +9.36%  '
$ws.Range("E10").Value = '  +2.94%  '
$ws.Range("D12").Formula = '="0.0994"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("D13").Formula = '="14.47"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  +9.27%  '
$ws.Range("E14").Value = '  +6.30%  '
$ws.Range("D15").Value = '2.197.61'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").Formula = '="5.12"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  +4.66%  '
$ws.Range("D17").Value = '1.913.04'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '36.449.43'
$ws.Range("E18").Value = '  +2.92%  '
$ws.Range("E19").Value = '  +1.96%  '
$ws.Range("E20").Value = '  +3.54%  '
$ws.Range("D21").Formula = '="252.17"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  +3.17%  '
$ws.Range("E22").Value = '  +3.26%  '
$ws.Range("D23").Formula = '="5.19"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +4.88%  '
$ws.Range("D24").Formula = '="2.71"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  +3.46%  '
$ws.Range("D27").Formula = '="167.79"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("E28").Value = '  +3.15%  '
$ws.Range("D29").Formula = '="18.85"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +3.09%  '
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("D31").Formula = '="4.52"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  +5.97%  '
$ws.Range("E32").Value = '  +4.05%  '
$ws.Range("E33").Value = '  +6.48%  '
$ws.Range("E34").Value = '  +4.48%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Formula = '="0.0853"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +22.91%  '
$ws.Range("E37").Value = '  -14.78%  '
$ws.Range("D38").Formula = '="0.861"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").Formula = '="2.00"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("D40").Formula = '="106.49"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +10.19%  '
$ws.Range("D41").Formula = '="16.28"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +32.64%  '
$ws.Range("E42").Value = '  +3.82%  '
$ws.Range("D43").Formula = '="16.99"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("E44").Value = '  +2.80%  '
$ws.Range("D45").Value = '1.341.63'
$ws.Range("E45").Value = '  +3.25%  '
$ws.Range("E46").Value = '  +1.89%  '
$ws.Range("D47").Formula = '="0.0807"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("E48").Value = '  +2.48%  '
$ws.Range("E49").Value = '  +2.20%  '
$ws.Range("E50").Value = '  +3.00%  '
$ws.Range("D51").Value = '2.097.58'
$ws.Range("E51").Value = '  +1.43%  '
